$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "491.71") as well as
# multi-dot thousands-grouped strings (e.g. "56.459.92"). Force text format on
# the whole price/volume column range first so assignment keeps the exact
# literal string instead of Excel coercing it to a floating-point number.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "56.459.92"
$ws.Range("E2").Value = "  +3.87%  "
$ws.Range("D3").Value = "2.511.08"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "491.71"
$ws.Range("E5").Value = "  +5.63%  "
$ws.Range("D6").Value = "147.24"
$ws.Range("E6").Value = "  +11.70%  "
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "0.516"
$ws.Range("E8").Value = "  +5.32%  "
$ws.Range("D9").Value = "2.529.67"
$ws.Range("E9").Value = "  +2.44%  "
$ws.Range("E10").Value = "  +7.40%  "
$ws.Range("D11").Value = "0.0980"
$ws.Range("E11").Value = "  +2.21%  "
$ws.Range("D12").Value = "0.335"
$ws.Range("E12").Value = "  +5.15%  "
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("D14").Value = "2.940.14"
$ws.Range("E14").Value = "  +2.18%  "
$ws.Range("D15").Value = "56.394.06"
$ws.Range("E15").Value = "  +4.17%  "
$ws.Range("D16").Value = "21.28"
$ws.Range("E16").Value = "  +7.47%  "
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D18").Value = "2.521.49"
$ws.Range("E18").Value = "  +1.99%  "
$ws.Range("D19").Value = "4.52"
$ws.Range("E19").Value = "  +7.52%  "
$ws.Range("D20").Value = "10.29"
$ws.Range("E20").Value = "  +8.73%  "
$ws.Range("D21").Value = "323.23"
$ws.Range("E21").Value = "  +3.51%  "
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").Value = "5.86"
$ws.Range("E23").Value = "  +8.99%  "
$ws.Range("D24").Value = "58.69"
$ws.Range("E24").Value = "  +3.61%  "
$ws.Range("D25").Value = "0.414"
$ws.Range("E25").Value = "  +7.58%  "
$ws.Range("D26").Value = "0.168"
$ws.Range("E26").Value = "  +8.40%  "
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("D28").Value = "2.616.67"
$ws.Range("E28").Value = "  +2.27%  "
$ws.Range("D29").Value = "7.63"
$ws.Range("E29").Value = "  +5.82%  "
$ws.Range("D30").Value = "0.0₃0806"
$ws.Range("E30").Value = "  +10.57%  "
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("D32").Value = "149.50"
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("E33").Value = "  +3.82%  "
$ws.Range("E34").Value = "  +5.65%  "
$ws.Range("D35").Value = "5.23"
$ws.Range("E35").Value = "  +4.15%  "
$ws.Range("D36").Value = "1.16"
$ws.Range("E36").Value = "  +8.93%  "
$ws.Range("D37").Value = "3.77"
$ws.Range("E37").Value = "  +5.33%  "
$ws.Range("D38").Value = "0.878"
$ws.Range("E38").Value = "  +10.84%  "
$ws.Range("D39").Value = "34.23"
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("D40").Value = "3.56"
$ws.Range("E40").Value = "  +8.63%  "
$ws.Range("D41").Value = "0.623"
$ws.Range("E41").Value = "  +3.09%  "
$ws.Range("D42").Value = "0.0559"
$ws.Range("E42").Value = "  +5.35%  "
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("E44").Value = "  +8.20%  "
$ws.Range("E45").Value = "  +11.98%  "
$ws.Range("D46").Value = "262.95"
$ws.Range("E46").Value = "  +18.32%  "
$ws.Range("D47").Value = "0.0231"
$ws.Range("E47").Value = "  +4.15%  "
$ws.Range("E48").Value = "  +5.15%  "
$ws.Range("D49").Value = "10.19"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").Value = "1.926.03"
$ws.Range("E50").Value = "  -1.96%  "
$ws.Range("D51").Value = "17.74"
$ws.Range("E51").Value = "  +5.81%  "
